$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 12:22"

# --- Update COVID stats (new data snapshot) ---

# Austria (row 20) - values updated, country order unchanged
$ws.Range("B20").Value = 14119
$ws.Range("C20").Value = 78
$ws.Range("E20").Value = 6102

# Rumania (row 31) - values updated, country order unchanged
$ws.Range("B31").Value = 6879
$ws.Range("C31").Value = 246
$ws.Range("D31").Value = 1051
$ws.Range("E31").Value = 5496
$ws.Range("F31").Value = 241

# Moldavia overtakes Marruecos (rows 59 & 60 swap order plus new values)
$ws.Range("A59").Value = "Moldavia"
$ws.Range("B59").Value = 1847
$ws.Range("C59").Value = 135
$ws.Range("D59").Value = 134
$ws.Range("E59").Value = 1677
$ws.Range("F59").Value = 80
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 36

$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 1838
$ws.Range("C60").Value = 75
$ws.Range("D60").Value = 210
$ws.Range("E60").Value = 1502
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 126

# Eslovenia (row 69) - values updated, country order unchanged
$ws.Range("B69").Value = 1220
$ws.Range("C69").Value = 8
$ws.Range("E69").Value = 1012
$ws.Range("F69").Value = 35
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 56

# Eslovaquia (row 80) - values updated, country order unchanged
$ws.Range("B80").Value = 835
$ws.Range("C80").Value = 66
$ws.Range("E80").Value = 726

# Vietnam (row 114) - values updated, country order unchanged
$ws.Range("D114").Value = 168
$ws.Range("E114").Value = 97

# Etiopia overtakes Liechtenstein and Togo (rows 137-139 shift)
$ws.Range("A137").Value = "Etiopia"
$ws.Range("B137").Value = 82
$ws.Range("C137").Value = 8
$ws.Range("D137").Value = 14
$ws.Range("E137").Value = 65
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 3

$ws.Range("A138").Value = "Liechtenstein"
$ws.Range("B138").Value = 79
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 55
$ws.Range("E138").Value = 23
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 1

$ws.Range("A139").Value = "Togo"
$ws.Range("B139").Value = 77
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 29
$ws.Range("E139").Value = 45
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 3

$wb.Save()
